$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (rows 2..101) held the full GitHub raw-image URL as the cell's
# displayed text (the hyperlink itself, pointing at the same URL, is a
# separate piece of metadata and is left untouched). The edit trims the
# display text down to just the bare filename, e.g.
#   https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A1.jpg?raw=true
# becomes
#   К1.jpg
# Image numbering goes 1..101 but skips 52 (there is no "К51"), matching
# row 52 holding "К52.jpg" straight after row 51's "К50.jpg".
$row = 2
for ($n = 1; $n -le 101; $n++) {
    if ($n -eq 51) { continue }
    $ws.Range("B$row").Value2 = "К$n.jpg"
    $row = $row + 1
}

# Update the current selection to match the saved view state.
$ws.Range("E108").Select()
